$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 7780.2  # ALC!H19: 8113.2666 -> 7780.2
$ws.Cells.Item(19, 10).Value = 7580.5  # ALC!J19: 8080.1 -> 7580.5
$ws.Cells.Item(19, 12).Value = 7580.5  # ALC!L19: 8080.1 -> 7580.5
$ws.Cells.Item(19, 14).Value = -7930.5  # ALC!N19: -8430.1 -> -7930.5
$ws.Cells.Item(33, 8).Value = 2565.1333  # ALC!H33: 2684.1428 -> 2565.1333
$ws.Cells.Item(33, 9).Value = 2614.8333  # ALC!I33: 2770.818 -> 2614.8333
$ws.Cells.Item(33, 11).Value = 2614.8333  # ALC!K33: 2770.818 -> 2614.8333
$ws.Cells.Item(33, 13).Value = -2385.8333  # ALC!M33: -2541.818 -> -2385.8333
$ws.Cells.Item(51, 8).Value = 14993.8  # ALC!H51: 14052.728 -> 14993.8
$ws.Cells.Item(51, 9).Value = 28966.666  # ALC!I51: 28832.666 -> 28966.666
$ws.Cells.Item(51, 10).Value = 9005.429  # ALC!J51: 8510.25 -> 9005.429
$ws.Cells.Item(51, 11).Value = 28966.666  # ALC!K51: 28832.666 -> 28966.666
$ws.Cells.Item(51, 12).Value = 9005.429  # ALC!L51: 8510.25 -> 9005.429
$ws.Cells.Item(51, 13).Value = -28482.666  # ALC!M51: -28348.666 -> -28482.666
$ws.Cells.Item(51, 14).Value = -9973.429  # ALC!N51: -9478.25 -> -9973.429
$ws.Cells.Item(87, 8).Value = 55999  # ALC!H87: 51999.25 -> 55999
$ws.Cells.Item(87, 10).Value = 55999  # ALC!J87: 51999.25 -> 55999
$ws.Cells.Item(87, 12).Value = 55999  # ALC!L87: 51999.25 -> 55999
$ws.Cells.Item(87, 14).Value = -58495  # ALC!N87: -54495.25 -> -58495
$ws.Cells.Item(90, 8).Value = 55999  # ALC!H90: 51999.25 -> 55999
$ws.Cells.Item(90, 10).Value = 55999  # ALC!J90: 51999.25 -> 55999
$ws.Cells.Item(90, 12).Value = 167997  # ALC!L90: 155997.75 -> 167997
$ws.Cells.Item(90, 14).Value = -180477  # ALC!N90: -168477.75 -> -180477
$ws.Cells.Item(98, 8).Value = 5495.6313  # ALC!H98: 6072.4707 -> 5495.6313
$ws.Cells.Item(98, 9).Value = 4915.8237  # ALC!I98: 5492.2666 -> 4915.8237
$ws.Cells.Item(98, 11).Value = 4915.8237  # ALC!K98: 5492.2666 -> 4915.8237
$ws.Cells.Item(98, 13).Value = -3417.8237  # ALC!M98: -3994.2666 -> -3417.8237
$ws.Cells.Item(121, 8).Value = 5350.7144  # ALC!H121: 6222.273 -> 5350.7144
$ws.Cells.Item(121, 9).Value = 2155  # ALC!I121: 0 -> 2155
$ws.Cells.Item(121, 11).Value = 6465  # ALC!K121: 0 -> 6465
$ws.Cells.Item(121, 13).Value = -4718  # ALC!M121: None -> -4718
$ws.Cells.Item(122, 8).Value = 5495.6313  # ALC!H122: 6072.4707 -> 5495.6313
$ws.Cells.Item(122, 9).Value = 4915.8237  # ALC!I122: 5492.2666 -> 4915.8237
$ws.Cells.Item(122, 11).Value = 14747.4711  # ALC!K122: 16476.7998 -> 14747.4711
$ws.Cells.Item(122, 13).Value = -12297.4711  # ALC!M122: -14026.7998 -> -12297.4711
$ws.Cells.Item(129, 8).Value = 1039.2941  # ALC!H129: 1041.6471 -> 1039.2941
$ws.Cells.Item(129, 9).Value = 550.0909  # ALC!I129: 553.7273 -> 550.0909
$ws.Cells.Item(129, 11).Value = 1650.2727  # ALC!K129: 1661.1819 -> 1650.2727
$ws.Cells.Item(129, 13).Value = 3349.7273  # ALC!M129: 3338.8181 -> 3349.7273
$ws.Cells.Item(132, 8).Value = 1795.8971  # ALC!H132: 1850.1538 -> 1795.8971
$ws.Cells.Item(132, 9).Value = 1730.0156  # ALC!I132: 1784.5902 -> 1730.0156
$ws.Cells.Item(132, 11).Value = 5190.0468  # ALC!K132: 5353.7706 -> 5190.0468
$ws.Cells.Item(132, 13).Value = -2660.0468  # ALC!M132: -2823.7706 -> -2660.0468
$ws.Cells.Item(137, 8).Value = 4515.121  # ALC!H137: 4606.2188 -> 4515.121
$ws.Cells.Item(137, 9).Value = 5966.6665  # ALC!I137: 6512.5 -> 5966.6665
$ws.Cells.Item(137, 11).Value = 17899.9995  # ALC!K137: 19537.5 -> 17899.9995
$ws.Cells.Item(137, 13).Value = -15349.9995  # ALC!M137: -16987.5 -> -15349.9995
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1766753.8  # ARM!H32: 1844694.4 -> 1766753.8
$ws.Cells.Item(32, 9).Value = 1898416.6  # ARM!I32: 1988812.8 -> 1898416.6
$ws.Cells.Item(32, 11).Value = 1898416.6  # ARM!K32: 1988812.8 -> 1898416.6
$ws.Cells.Item(32, 13).Value = -1898129.6  # ARM!M32: -1988525.8 -> -1898129.6
$ws.Cells.Item(45, 8).Value = 5904.5  # ARM!H45: 4358 -> 5904.5
$ws.Cells.Item(45, 9).Value = 0  # ARM!I45: 1265 -> 0
$ws.Cells.Item(45, 11).Value = 0  # ARM!K45: 1265 -> 0
$ws.Cells.Item(45, 13).ClearContents()  # ARM!M45: remove (was -888)
$ws.Cells.Item(61, 8).Value = 22226242  # ARM!H61: 24394596 -> 22226242
$ws.Cells.Item(61, 9).Value = 3078.6487  # ARM!I61: 3376.2424 -> 3078.6487
$ws.Cells.Item(61, 11).Value = 3078.6487  # ARM!K61: 3376.2424 -> 3078.6487
$ws.Cells.Item(61, 13).Value = -2866.6487  # ARM!M61: -3164.2424 -> -2866.6487
$ws.Cells.Item(74, 8).Value = 2588.7778  # ARM!H74: 2700.1372 -> 2588.7778
$ws.Cells.Item(74, 9).Value = 1856.9524  # ARM!I74: 1946.2821 -> 1856.9524
$ws.Cells.Item(74, 11).Value = 1856.9524  # ARM!K74: 1946.2821 -> 1856.9524
$ws.Cells.Item(74, 13).Value = -982.9523999999999  # ARM!M74: -1072.2821 -> -982.9523999999999
$ws.Cells.Item(77, 8).Value = 2588.7778  # ARM!H77: 2700.1372 -> 2588.7778
$ws.Cells.Item(77, 9).Value = 1856.9524  # ARM!I77: 1946.2821 -> 1856.9524
$ws.Cells.Item(77, 11).Value = 9284.761999999999  # ARM!K77: 9731.4105 -> 9284.761999999999
$ws.Cells.Item(77, 13).Value = -4916.761999999999  # ARM!M77: -5363.4105 -> -4916.761999999999
$ws.Cells.Item(122, 8).Value = 4310.4287  # ARM!H122: 4295.0835 -> 4310.4287
$ws.Cells.Item(122, 9).Value = 3237.375  # ARM!I122: 3258.2 -> 3237.375
$ws.Cells.Item(122, 11).Value = 9712.125  # ARM!K122: 9774.599999999999 -> 9712.125
$ws.Cells.Item(122, 13).Value = -7262.125  # ARM!M122: -7324.599999999999 -> -7262.125
$ws.Cells.Item(136, 8).Value = 22226242  # ARM!H136: 24394596 -> 22226242
$ws.Cells.Item(136, 9).Value = 3078.6487  # ARM!I136: 3376.2424 -> 3078.6487
$ws.Cells.Item(136, 11).Value = 9235.946100000001  # ARM!K136: 10128.7272 -> 9235.946100000001
$ws.Cells.Item(136, 13).Value = -6685.946100000001  # ARM!M136: -7578.727200000001 -> -6685.946100000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 4764565  # BSM!H20: 4764615.5 -> 4764565
$ws.Cells.Item(20, 9).Value = 9261580  # BSM!I20: 9261678 -> 9261580
$ws.Cells.Item(20, 11).Value = 9261580  # BSM!K20: 9261678 -> 9261580
$ws.Cells.Item(20, 13).Value = -9261333  # BSM!M20: -9261431 -> -9261333
$ws.Cells.Item(22, 8).Value = 330.77777  # BSM!H22: 312.7 -> 330.77777
$ws.Cells.Item(22, 10).Value = 475  # BSM!J22: 410 -> 475
$ws.Cells.Item(22, 12).Value = 475  # BSM!L22: 410 -> 475
$ws.Cells.Item(22, 14).Value = -821  # BSM!N22: -756 -> -821
$ws.Cells.Item(99, 8).Value = 8266774  # BSM!H99: 9093262 -> 8266774
$ws.Cells.Item(99, 9).Value = 2419.1428  # BSM!I99: 2433.4285 -> 2419.1428
$ws.Cells.Item(99, 10).Value = 22729396  # BSM!J99: 30305194 -> 22729396
$ws.Cells.Item(99, 11).Value = 2419.1428  # BSM!K99: 2433.4285 -> 2419.1428
$ws.Cells.Item(99, 12).Value = 22729396  # BSM!L99: 30305194 -> 22729396
$ws.Cells.Item(99, 13).Value = -921.1428000000001  # BSM!M99: -935.4285 -> -921.1428000000001
$ws.Cells.Item(99, 14).Value = -22732392  # BSM!N99: -30308190 -> -22732392
$ws.Cells.Item(105, 8).Value = 3991.6667  # BSM!H105: 4080.3704 -> 3991.6667
$ws.Cells.Item(105, 9).Value = 3194.25  # BSM!I105: 3286.6155 -> 3194.25
$ws.Cells.Item(105, 10).Value = 4629.6  # BSM!J105: 4817.4287 -> 4629.6
$ws.Cells.Item(105, 11).Value = 3194.25  # BSM!K105: 3286.6155 -> 3194.25
$ws.Cells.Item(105, 12).Value = 4629.6  # BSM!L105: 4817.4287 -> 4629.6
$ws.Cells.Item(105, 13).Value = -1447.25  # BSM!M105: -1539.6155 -> -1447.25
$ws.Cells.Item(105, 14).Value = -8123.6  # BSM!N105: -8311.4287 -> -8123.6
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 8976.947  # CRP!H31: 9162.973 -> 8976.947
$ws.Cells.Item(31, 9).Value = 4008.7646  # CRP!I31: 4128.4375 -> 4008.7646
$ws.Cells.Item(31, 11).Value = 4008.7646  # CRP!K31: 4128.4375 -> 4008.7646
$ws.Cells.Item(31, 13).Value = -3713.7646  # CRP!M31: -3833.4375 -> -3713.7646
$ws.Cells.Item(34, 8).Value = 8976.947  # CRP!H34: 9162.973 -> 8976.947
$ws.Cells.Item(34, 9).Value = 4008.7646  # CRP!I34: 4128.4375 -> 4008.7646
$ws.Cells.Item(34, 11).Value = 4008.7646  # CRP!K34: 4128.4375 -> 4008.7646
$ws.Cells.Item(34, 13).Value = -3806.7646  # CRP!M34: -3926.4375 -> -3806.7646
$ws.Cells.Item(58, 8).Value = 6800.049  # CRP!H58: 7278.1313 -> 6800.049
$ws.Cells.Item(58, 9).Value = 2451.5  # CRP!I58: 2845.4614 -> 2451.5
$ws.Cells.Item(58, 11).Value = 2451.5  # CRP!K58: 2845.4614 -> 2451.5
$ws.Cells.Item(58, 13).Value = -2248.5  # CRP!M58: -2642.4614 -> -2248.5
$ws.Cells.Item(99, 8).Value = 6346.905  # CRP!H99: 6348.8096 -> 6346.905
$ws.Cells.Item(99, 9).Value = 4659.3335  # CRP!I99: 4659.5 -> 4659.3335
$ws.Cells.Item(99, 10).Value = 6628.1665  # CRP!J99: 6526.6313 -> 6628.1665
$ws.Cells.Item(99, 11).Value = 4659.3335  # CRP!K99: 4659.5 -> 4659.3335
$ws.Cells.Item(99, 12).Value = 6628.1665  # CRP!L99: 6526.6313 -> 6628.1665
$ws.Cells.Item(99, 13).Value = -3161.3335  # CRP!M99: -3161.5 -> -3161.3335
$ws.Cells.Item(99, 14).Value = -9624.166499999999  # CRP!N99: -9522.631300000001 -> -9624.166499999999
$ws.Cells.Item(126, 8).Value = 6346.905  # CRP!H126: 6348.8096 -> 6346.905
$ws.Cells.Item(126, 9).Value = 4659.3335  # CRP!I126: 4659.5 -> 4659.3335
$ws.Cells.Item(126, 10).Value = 6628.1665  # CRP!J126: 6526.6313 -> 6628.1665
$ws.Cells.Item(126, 11).Value = 13978.0005  # CRP!K126: 13978.5 -> 13978.0005
$ws.Cells.Item(126, 12).Value = 19884.4995  # CRP!L126: 19579.8939 -> 19884.4995
$ws.Cells.Item(126, 13).Value = -11508.0005  # CRP!M126: -11508.5 -> -11508.0005
$ws.Cells.Item(126, 14).Value = -24824.4995  # CRP!N126: -24519.8939 -> -24824.4995
$ws.Cells.Item(134, 8).Value = 5471.7114  # CRP!H134: 5791.2915 -> 5471.7114
$ws.Cells.Item(134, 9).Value = 2001.1936  # CRP!I134: 2160.2856 -> 2001.1936
$ws.Cells.Item(134, 10).Value = 10594.857  # CRP!J134: 10874.7 -> 10594.857
$ws.Cells.Item(134, 11).Value = 6003.5808  # CRP!K134: 6480.8568 -> 6003.5808
$ws.Cells.Item(134, 12).Value = 31784.571  # CRP!L134: 32624.1 -> 31784.571
$ws.Cells.Item(134, 13).Value = -3468.5808  # CRP!M134: -3945.8568 -> -3468.5808
$ws.Cells.Item(134, 14).Value = -36854.571  # CRP!N134: -37694.10000000001 -> -36854.571
$ws.Cells.Item(136, 8).Value = 6800.049  # CRP!H136: 7278.1313 -> 6800.049
$ws.Cells.Item(136, 9).Value = 2451.5  # CRP!I136: 2845.4614 -> 2451.5
$ws.Cells.Item(136, 11).Value = 7354.5  # CRP!K136: 8536.3842 -> 7354.5
$ws.Cells.Item(136, 13).Value = -4804.5  # CRP!M136: -5986.3842 -> -4804.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 742.4048  # GSM!H97: 731.3022999999999 -> 742.4048
$ws.Cells.Item(97, 9).Value = 593.7692  # GSM!I97: 581.5925999999999 -> 593.7692
$ws.Cells.Item(97, 11).Value = 593.7692  # GSM!K97: 581.5925999999999 -> 593.7692
$ws.Cells.Item(97, 13).Value = -97.76919999999996  # GSM!M97: -85.59259999999995 -> -97.76919999999996
$ws.Cells.Item(102, 8).Value = 2959.0588  # GSM!H102: 2966.6667 -> 2959.0588
$ws.Cells.Item(102, 9).Value = 2959.0588  # GSM!I102: 2966.6667 -> 2959.0588
$ws.Cells.Item(102, 11).Value = 2959.0588  # GSM!K102: 2966.6667 -> 2959.0588
$ws.Cells.Item(102, 13).Value = -1337.0588  # GSM!M102: -1344.6667 -> -1337.0588
$ws.Cells.Item(132, 8).Value = 2542.0715  # GSM!H132: 2481.6553 -> 2542.0715
$ws.Cells.Item(132, 9).Value = 1444.5294  # GSM!I132: 1443.9412 -> 1444.5294
$ws.Cells.Item(132, 10).Value = 4238.273  # GSM!J132: 3951.75 -> 4238.273
$ws.Cells.Item(132, 11).Value = 4333.5882  # GSM!K132: 4331.8236 -> 4333.5882
$ws.Cells.Item(132, 12).Value = 12714.819  # GSM!L132: 11855.25 -> 12714.819
$ws.Cells.Item(132, 13).Value = -1803.5882  # GSM!M132: -1801.8236 -> -1803.5882
$ws.Cells.Item(132, 14).Value = -17774.819  # GSM!N132: -16915.25 -> -17774.819
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 749.05554  # LTW!H16: 755.5294 -> 749.05554
$ws.Cells.Item(16, 9).Value = 749.7143  # LTW!I16: 758.2308 -> 749.7143
$ws.Cells.Item(16, 11).Value = 749.7143  # LTW!K16: 758.2308 -> 749.7143
$ws.Cells.Item(16, 13).Value = -579.7143  # LTW!M16: -588.2308 -> -579.7143
$ws.Cells.Item(22, 8).Value = 2033  # LTW!H22: 1892 -> 2033
$ws.Cells.Item(22, 9).Value = 711  # LTW!I22: 698.6667 -> 711
$ws.Cells.Item(22, 10).Value = 3355  # LTW!J22: 3085.3333 -> 3355
$ws.Cells.Item(22, 11).Value = 711  # LTW!K22: 698.6667 -> 711
$ws.Cells.Item(22, 12).Value = 3355  # LTW!L22: 3085.3333 -> 3355
$ws.Cells.Item(22, 13).Value = -416  # LTW!M22: -403.6667 -> -416
$ws.Cells.Item(22, 14).Value = -3945  # LTW!N22: -3675.3333 -> -3945
$ws.Cells.Item(27, 8).Value = 2033  # LTW!H27: 1892 -> 2033
$ws.Cells.Item(27, 9).Value = 711  # LTW!I27: 698.6667 -> 711
$ws.Cells.Item(27, 10).Value = 3355  # LTW!J27: 3085.3333 -> 3355
$ws.Cells.Item(27, 11).Value = 711  # LTW!K27: 698.6667 -> 711
$ws.Cells.Item(27, 12).Value = 3355  # LTW!L27: 3085.3333 -> 3355
$ws.Cells.Item(27, 13).Value = -604  # LTW!M27: -591.6667 -> -604
$ws.Cells.Item(27, 14).Value = -3569  # LTW!N27: -3299.3333 -> -3569
$ws.Cells.Item(61, 8).Value = 5339.0454  # LTW!H61: 6034.737 -> 5339.0454
$ws.Cells.Item(61, 9).Value = 3138.7144  # LTW!I61: 3528.5 -> 3138.7144
$ws.Cells.Item(61, 10).Value = 6365.8667  # LTW!J61: 7191.4614 -> 6365.8667
$ws.Cells.Item(61, 11).Value = 3138.7144  # LTW!K61: 3528.5 -> 3138.7144
$ws.Cells.Item(61, 12).Value = 6365.8667  # LTW!L61: 7191.4614 -> 6365.8667
$ws.Cells.Item(61, 13).Value = -2936.7144  # LTW!M61: -3326.5 -> -2936.7144
$ws.Cells.Item(61, 14).Value = -6769.8667  # LTW!N61: -7595.4614 -> -6769.8667
$ws.Cells.Item(100, 8).Value = 5162.875  # LTW!H100: 5200.4287 -> 5162.875
$ws.Cells.Item(100, 10).Value = 5515.1816  # LTW!J100: 5651.8887 -> 5515.1816
$ws.Cells.Item(100, 12).Value = 5515.1816  # LTW!L100: 5651.8887 -> 5515.1816
$ws.Cells.Item(100, 14).Value = -6597.1816  # LTW!N100: -6733.8887 -> -6597.1816
$ws.Cells.Item(113, 8).Value = 5339.0454  # LTW!H113: 6034.737 -> 5339.0454
$ws.Cells.Item(113, 9).Value = 3138.7144  # LTW!I113: 3528.5 -> 3138.7144
$ws.Cells.Item(113, 10).Value = 6365.8667  # LTW!J113: 7191.4614 -> 6365.8667
$ws.Cells.Item(113, 11).Value = 3138.7144  # LTW!K113: 3528.5 -> 3138.7144
$ws.Cells.Item(113, 12).Value = 6365.8667  # LTW!L113: 7191.4614 -> 6365.8667
$ws.Cells.Item(113, 13).Value = -968.7143999999998  # LTW!M113: -1358.5 -> -968.7143999999998
$ws.Cells.Item(113, 14).Value = -10705.8667  # LTW!N113: -11531.4614 -> -10705.8667
$ws.Cells.Item(132, 8).Value = 7580775.5  # LTW!H132: 8069774.5 -> 7580775.5
$ws.Cells.Item(132, 9).Value = 12197354  # LTW!I132: 13515846 -> 12197354
$ws.Cells.Item(132, 11).Value = 36592062  # LTW!K132: 40547538 -> 36592062
$ws.Cells.Item(132, 13).Value = -36589532  # LTW!M132: -40545008 -> -36589532
$ws.Cells.Item(136, 8).Value = 7641.551  # LTW!H136: 7771.213 -> 7641.551
$ws.Cells.Item(136, 9).Value = 2622.84  # LTW!I136: 2673.6667 -> 2622.84
$ws.Cells.Item(136, 10).Value = 12869.375  # LTW!J136: 13090.392 -> 12869.375
$ws.Cells.Item(136, 11).Value = 7868.52  # LTW!K136: 8021.000100000001 -> 7868.52
$ws.Cells.Item(136, 12).Value = 38608.125  # LTW!L136: 39271.176 -> 38608.125
$ws.Cells.Item(136, 13).Value = -5318.52  # LTW!M136: -5471.000100000001 -> -5318.52
$ws.Cells.Item(136, 14).Value = -43708.125  # LTW!N136: -44371.176 -> -43708.125
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(26, 8).Value = 0  # WVR!H26: 10000 -> 0
$ws.Cells.Item(26, 10).Value = 0  # WVR!J26: 10000 -> 0
$ws.Cells.Item(26, 12).Value = 0  # WVR!L26: 10000 -> 0
$ws.Cells.Item(26, 14).ClearContents()  # WVR!N26: remove (was -10586)
$ws.Cells.Item(41, 8).Value = 11634.454  # WVR!H41: 12091.363 -> 11634.454
$ws.Cells.Item(41, 10).Value = 10069.4  # WVR!J41: 10572 -> 10069.4
$ws.Cells.Item(41, 12).Value = 10069.4  # WVR!L41: 10572 -> 10069.4
$ws.Cells.Item(41, 14).Value = -10849.4  # WVR!N41: -11352 -> -10849.4
$ws.Cells.Item(113, 8).Value = 1007.88464  # WVR!H113: 1036.5 -> 1007.88464
$ws.Cells.Item(113, 9).Value = 847.7619  # WVR!I113: 871.15 -> 847.7619
$ws.Cells.Item(113, 10).Value = 1680.4  # WVR!J113: 1863.25 -> 1680.4
$ws.Cells.Item(113, 11).Value = 2543.2857  # WVR!K113: 2613.45 -> 2543.2857
$ws.Cells.Item(113, 12).Value = 5041.200000000001  # WVR!L113: 5589.75 -> 5041.200000000001
$ws.Cells.Item(113, 13).Value = -373.2856999999999  # WVR!M113: -443.4499999999998 -> -373.2856999999999
$ws.Cells.Item(113, 14).Value = -9381.200000000001  # WVR!N113: -9929.75 -> -9381.200000000001
$ws.Cells.Item(122, 8).Value = 100656.93  # WVR!H122: 91782.78 -> 100656.93
$ws.Cells.Item(122, 9).Value = 139984.14  # WVR!I122: 123116.09 -> 139984.14
$ws.Cells.Item(122, 11).Value = 419952.42  # WVR!K122: 369348.27 -> 419952.42
$ws.Cells.Item(122, 13).Value = -417502.42  # WVR!M122: -366898.27 -> -417502.42
$ws.Cells.Item(136, 8).Value = 17721284  # WVR!H136: 18037700 -> 17721284
$ws.Cells.Item(136, 10).Value = 672200.4  # WVR!J136: 720072 -> 672200.4
$ws.Cells.Item(136, 12).Value = 2016601.2  # WVR!L136: 2160216 -> 2016601.2
$ws.Cells.Item(136, 14).Value = -2021701.2  # WVR!N136: -2165316 -> -2021701.2
